$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the explicit 28.5pt row height from rows 2-16 (revert to the sheet's
# default row height) by auto-fitting them.
$ws.Rows("2:16").AutoFit()

# Widen column A (was 18.265625 chars ~ "27.06640625" stored width). The
# ColumnWidth setter is expressed in "characters"; feed it the value whose
# round-trip lands closest to the target stored width.
$ws.Columns("A").ColumnWidth = 26.166666666666668

# Move the selection to B20.
$ws.Range("B20").Select() | Out-Null

# Zoom the view to 150%.
$excel.ActiveWindow.Zoom = 150
